$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking
# strings like "1.000" / "29.323.81" are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.323.81"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.876.58"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("D5").Value = "0.7109"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "241.98"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.08103"
$ws.Range("E8").Value = "  +4.86%  "
$ws.Range("D9").Value = "0.3124"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "0.08391"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "1.877.29"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "5.247"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").Value = "0.7182"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "91.30"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "6.223"
$ws.Range("E16").Value = "  +3.93%  "
$ws.Range("D17").Value = "0.000008400"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "29.319.47"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "240.72"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "13.24"
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("D21").Value = "2.126.28"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "7.793"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "0.9991"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "0.1594"
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("D26").Value = "162.84"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "9.056"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").Value = "18.56"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "4.423"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").Value = "4.349"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("D32").Value = "1.206"
$ws.Range("E32").Value = "  -6.25%  "
$ws.Range("D33").Value = "0.05359"
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("D34").Value = "1.946"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.177"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.7501"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").Value = "2.696"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").Value = "1.288.85"
$ws.Range("E38").Value = "  +11.56%  "
$ws.Range("D39").Value = "0.01886"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("D40").Value = "2.736"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("E41").Value = "  +3.41%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "110.97"
$ws.Range("E42").Value = "  +4.49%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.8935"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("D44").Value = "73.31"
$ws.Range("E45").Value = "  +9.00%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "2.019.91"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "1.800"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").Value = "0.5202"
$ws.Range("D50").Value = "9.462"
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("D51").Value = "0.4365"
$ws.Range("E51").Value = "  +1.50%  "
